# Correction des fautes d'orthographe (Scénarios_teste.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Header: Arthur -> "Arthur" + newline + "HP ENVY" -----------------
# Setting the cell value directly keeps the ListObject's header in sync
# (table1.xml tableColumn name is driven from the header cell).
$ws.Range("D4").Value = "Arthur" + [char]10 + "HP ENVY"

# --- Spelling / wording corrections in column B ------------------------
$ws.Range("B6").Value  = "1.1 Taper 2 et ensuite enter"
$ws.Range("B7").Value  = "1.2 Lire les règles et appuyer sur Enter"
$ws.Range("B10").Value = "2.3 Taper sur N et ensuite  Enter"
$ws.Range("B15").Value = "3.3 Entrer les coordonnées 838-AJF"
$ws.Range("B18").Value = "3.6 Un écran s'affiche, vous félicite et vous donne votre score"
$ws.Range("B20").Value = "4.1 Appuyez sur 5 et ensuite Enter"
$ws.Range("B21").Value = "4.2 Vérifier que votre nom et votre score soient affichés et appuyer sur enter"
$ws.Range("B26").Value = "6.3 Le programme se ferme"

# --- Layout tweaks -------------------------------------------------------
# Column B: narrower, best-fit width instead of the old fixed width
# (~70.29 chars once text shrank after the spelling clean-up / AutoFit).
$ws.Columns.Item(2).ColumnWidth = 69.5

# Row 4 (header row) grows to fit the two-line "Arthur / HP ENVY" header.
$ws.Rows.Item(4).RowHeight = 30

# Reset the remembered selection back to the top-left cell (no stray
# "G20" selection saved with the workbook).
$ws.Range("A1").Select()
